$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 22.700661
$ws.Range("H2").Value = 68.10198299999999
$ws.Range("I2").Value = 0.08615268874617349
$ws.Range("J2").Value = 0.08615268874617349
$ws.Range("M2").Value = 1.168007333333333
$ws.Range("N2").Value = 3.504022
$ws.Range("O2").Value = 0.1638609704511517
$ws.Range("P2").Value = 0.1638609704511517
$ws.Range("Q2").Value = 26.514538519514
$ws.Range("R2").Value = 238.630846675626
$ws.Range("S2").Value = 0.01411706318492401
$ws.Range("T2").Value = 0.01411706318492401
$ws.Range("G3").Value = 22.700661
$ws.Range("H3").Value = 68.10198299999999
$ws.Range("I3").Value = 0.08615268874617349
$ws.Range("J3").Value = 0.08615268874617349
$ws.Range("O3").Value = 0.5019752511630595
$ws.Range("P3").Value = 0.5019752511630595
$ws.Range("Q3").Value = 81.22521242343898
$ws.Range("R3").Value = 731.0269118109508
$ws.Range("S3").Value = 0.04324651757173333
$ws.Range("T3").Value = 0.04324651757173333
$ws.Range("G4").Value = 22.700661
$ws.Range("H4").Value = 68.10198299999999
$ws.Range("I4").Value = 0.08615268874617349
$ws.Range("J4").Value = 0.08615268874617349
$ws.Range("M4").Value = 1.915392333333333
$ws.Range("N4").Value = 5.746177
$ws.Range("O4").Value = 0.2687123938160456
$ws.Range("P4").Value = 0.2687123938160456
$ws.Range("Q4").Value = 43.480672040999
$ws.Range("R4").Value = 391.3260483689909
$ws.Range("S4").Value = 0.02315029522667297
$ws.Range("T4").Value = 0.02315029522667297
$ws.Range("G5").Value = 22.700661
$ws.Range("H5").Value = 68.10198299999999
$ws.Range("I5").Value = 0.08615268874617349
$ws.Range("J5").Value = 0.08615268874617349
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.46654
$ws.Range("N5").Value = 1.39962
$ws.Range("O5").Value = 0.06545138456974327
$ws.Range("P5").Value = 0.06545138456974327
$ws.Range("Q5").Value = 10.59076638294
$ws.Range("R5").Value = 95.31689744645998
$ws.Range("S5").Value = 0.005638812762843194
$ws.Range("T5").Value = 0.005638812762843194
$ws.Range("I6").Value = 0.5030288587986086
$ws.Range("J6").Value = 0.5030288587986087
$ws.Range("M6").Value = 1.168007333333333
$ws.Range("N6").Value = 3.504022
$ws.Range("O6").Value = 0.1638609704511517
$ws.Range("P6").Value = 0.1638609704511517
$ws.Range("Q6").Value = 154.8132536215855
$ws.Range("R6").Value = 1393.31928259427
$ws.Range("S6").Value = 0.08242679696767538
$ws.Range("T6").Value = 0.08242679696767539
$ws.Range("I7").Value = 0.5030288587986086
$ws.Range("J7").Value = 0.5030288587986087
$ws.Range("O7").Value = 0.5019752511630595
$ws.Range("P7").Value = 0.5019752511630595
$ws.Range("S7").Value = 0.2525080377376988
$ws.Range("T7").Value = 0.2525080377376988
$ws.Range("I8").Value = 0.5030288587986086
$ws.Range("J8").Value = 0.5030288587986087
$ws.Range("M8").Value = 1.915392333333333
$ws.Range("N8").Value = 5.746177
$ws.Range("O8").Value = 0.2687123938160456
$ws.Range("P8").Value = 0.2687123938160456
$ws.Range("Q8").Value = 253.8752203198272
$ws.Range("R8").Value = 2284.876982878445
$ws.Range("S8").Value = 0.1351700888063277
$ws.Range("T8").Value = 0.1351700888063277
$ws.Range("I9").Value = 0.5030288587986086
$ws.Range("J9").Value = 0.5030288587986087
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.46654
$ws.Range("N9").Value = 1.39962
$ws.Range("O9").Value = 0.06545138456974327
$ws.Range("P9").Value = 0.06545138456974327
$ws.Range("Q9").Value = 61.83743310796666
$ws.Range("R9").Value = 556.5368979716999
$ws.Range("S9").Value = 0.03292393528690681
$ws.Range("T9").Value = 0.03292393528690683
$ws.Range("G10").Value = 41.94534433333333
$ws.Range("H10").Value = 125.836033
$ws.Range("I10").Value = 0.159189381961201
$ws.Range("J10").Value = 0.159189381961201
$ws.Range("M10").Value = 1.168007333333333
$ws.Range("N10").Value = 3.504022
$ws.Range("O10").Value = 0.1638609704511517
$ws.Range("P10").Value = 0.1638609704511517
$ws.Range("Q10").Value = 48.99246978052511
$ws.Range("R10").Value = 440.9322280247259
$ws.Range("S10").Value = 0.02608492661368146
$ws.Range("T10").Value = 0.02608492661368146
$ws.Range("G11").Value = 41.94534433333333
$ws.Range("H11").Value = 125.836033
$ws.Range("I11").Value = 0.159189381961201
$ws.Range("J11").Value = 0.159189381961201
$ws.Range("O11").Value = 0.5019752511630595
$ws.Range("P11").Value = 0.5019752511630595
$ws.Range("Q11").Value = 150.0845946137556
$ws.Range("R11").Value = 1350.761351523801
$ws.Range("S11").Value = 0.07990912999246609
$ws.Range("T11").Value = 0.07990912999246609
$ws.Range("G12").Value = 41.94534433333333
$ws.Range("H12").Value = 125.836033
$ws.Range("I12").Value = 0.159189381961201
$ws.Range("J12").Value = 0.159189381961201
$ws.Range("M12").Value = 1.915392333333333
$ws.Range("N12").Value = 5.746177
$ws.Range("O12").Value = 0.2687123938160456
$ws.Range("P12").Value = 0.2687123938160456
$ws.Range("Q12").Value = 80.34179095509344
$ws.Range("R12").Value = 723.076118595841
$ws.Range("S12").Value = 0.04277615989689114
$ws.Range("T12").Value = 0.04277615989689115
$ws.Range("G13").Value = 41.94534433333333
$ws.Range("H13").Value = 125.836033
$ws.Range("I13").Value = 0.159189381961201
$ws.Range("J13").Value = 0.159189381961201
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.46654
$ws.Range("N13").Value = 1.39962
$ws.Range("O13").Value = 0.06545138456974327
$ws.Range("P13").Value = 0.06545138456974327
$ws.Range("Q13").Value = 19.56918094527333
$ws.Range("R13").Value = 176.12262850746
$ws.Range("S13").Value = 0.01041916545816232
$ws.Range("T13").Value = 0.01041916545816232
$ws.Range("G14").Value = 66.302588
$ws.Range("H14").Value = 198.907764
$ws.Range("I14").Value = 0.2516290704940168
$ws.Range("J14").Value = 0.2516290704940168
$ws.Range("M14").Value = 1.168007333333333
$ws.Range("N14").Value = 3.504022
$ws.Range("O14").Value = 0.1638609704511517
$ws.Range("P14").Value = 0.1638609704511517
$ws.Range("Q14").Value = 77.44190900297868
$ws.Range("R14").Value = 696.9771810268079
$ws.Range("S14").Value = 0.04123218368487087
$ws.Range("T14").Value = 0.04123218368487087
$ws.Range("G15").Value = 66.302588
$ws.Range("H15").Value = 198.907764
$ws.Range("I15").Value = 0.2516290704940168
$ws.Range("J15").Value = 0.2516290704940168
$ws.Range("O15").Value = 0.5019752511630595
$ws.Range("P15").Value = 0.5019752511630595
$ws.Range("Q15").Value = 237.237223820212
$ws.Range("R15").Value = 2135.135014381907
$ws.Range("S15").Value = 0.1263115658611613
$ws.Range("T15").Value = 0.1263115658611613
$ws.Range("G16").Value = 66.302588
$ws.Range("H16").Value = 198.907764
$ws.Range("I16").Value = 0.2516290704940168
$ws.Range("J16").Value = 0.2516290704940168
$ws.Range("M16").Value = 1.915392333333333
$ws.Range("N16").Value = 5.746177
$ws.Range("O16").Value = 0.2687123938160456
$ws.Range("P16").Value = 0.2687123938160456
$ws.Range("Q16").Value = 126.9954687353587
$ws.Range("R16").Value = 1142.959218618228
$ws.Range("S16").Value = 0.06761584988615374
$ws.Range("T16").Value = 0.06761584988615375
$ws.Range("G17").Value = 66.302588
$ws.Range("H17").Value = 198.907764
$ws.Range("I17").Value = 0.2516290704940168
$ws.Range("J17").Value = 0.2516290704940168
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.46654
$ws.Range("N17").Value = 1.39962
$ws.Range("O17").Value = 0.06545138456974327
$ws.Range("P17").Value = 0.06545138456974327
$ws.Range("Q17").Value = 30.93280940552
$ws.Range("R17").Value = 278.39528464968
$ws.Range("S17").Value = 0.01646947106183094
$ws.Range("T17").Value = 0.01646947106183094
